$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.551.80"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "1.563.86"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "'0.989"
$ws.Range("E4").Value = "  -1.71%  "
$ws.Range("D5").Value = "'210.77"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  -1.77%  "
$ws.Range("D8").Value = "'22.62"
$ws.Range("E8").Value = "  +2.41%  "
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.563.52"
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D16").Value = "27.507.71"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").Value = "'62.50"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "'225.01"
$ws.Range("E18").Value = "  +4.35%  "
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").Value = "0.0₃0705"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "'9.41"
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "'149.96"
$ws.Range("E25").Value = "  -2.23%  "
# Row 26: Stellar -> EthereumClassic (swap with row 28)
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'15.19"
$ws.Range("E26").Value = "  +1.04%  "

$ws.Range("E27").Value = "  +0.35%  "
# Row 28: EthereumClassic -> Stellar (swap with row 26)
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "'0.108"
$ws.Range("E28").Value = "  +2.26%  "

$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "1.455.62"
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("D34").Value = "'3.17"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").Value = "'0.0168"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "'0.542"
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("D40").Value = "'0.815"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").Value = "'0.988"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").Value = "'1.85"
$ws.Range("E44").Value = "  +6.60%  "
$ws.Range("D45").Value = "'0.969"
$ws.Range("E45").Value = "  -3.63%  "
$ws.Range("D46").Value = "'64.77"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "1.700.07"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "'86.45"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("D49").Value = "'0.0524"
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("D50").Value = "0.0₇0973"
$ws.Range("E50").Value = "  -6.28%  "
$ws.Range("E51").Value = "  -1.12%  "
